$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for d77931b3...md row
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-04 23:09:40"

# zh-cn sheet: "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# for d77931b3...md row
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-04 23:09:35"
$wsZhCn.Range("K2").Value = "2016-09-04 23:10:14"

# de-de sheet: "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# for d77931b3...md row
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-09-04 23:09:40"
$wsDeDe.Range("K2").Value = "2016-09-04 23:10:23"
